# feat(CWL): `SafeValueBonus` patch for `ElementContainer`
#
# Updates six localization strings in the "General" sheet: each of these
# rows keeps identical text in columns C (text_JP) and D (text), so both
# get the same new value. Excel will re-wrap the shared-string table and
# auto-recompute the (non customHeight) row heights for the wrapped cells
# on its own — no manual index/height bookkeeping required here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cwl_warn_deserialize (row 46): drop the "CWL caught the exception..." line
$newDeserialize = "failed to create {0} id: {1}, type: {2}, it may be missing from current game`nif this is causing issues,  please check for outdated mods or disable {3} in the config file"
$ws.Range("C46").Value = $newDeserialize
$ws.Range("D46").Value = $newDeserialize

# cwl_type_safety_desc (row 49): drop the "CWL kept the game going..." line
# and the blank line that followed it
$newSafetyDesc = "This class is missing or modified from your current game. `nYou should report the relevant information to mod author or CWL.`nUsing this element(if usable) will let CWL purge it from your save.`nYou may also keep the safety cone, CWL will restore it when the responsible mod functions again."
$ws.Range("C49").Value = $newSafetyDesc
$ws.Range("D49").Value = $newSafetyDesc

# cwl_warn_quest_id_thing (row 80): reword second line
$newQuestId = "quest {0} is trying to use invalid id: `"{1}`"`nCWL replaced it with `"{2}`""
$ws.Range("C80").Value = $newQuestId
$ws.Range("D80").Value = $newQuestId

# cwl_warn_drama_play_ex (row 112): drop the "CWL caught the exception..." line
$newDramaPlayEx = "Error occurred during drama play!`nPlease check the Player.log and mods.`n{0}"
$ws.Range("C112").Value = $newDramaPlayEx
$ws.Range("D112").Value = $newDramaPlayEx

# cwl_warn_fix_actCombat (row 120): prefix with "CWL "
$newActCombat = "CWL removed invalid actCombat ID: {0} from {1}"
$ws.Range("C120").Value = $newActCombat
$ws.Range("D120").Value = $newActCombat

# cwl_warn_fix_listAbility (row 121): prefix with "CWL "
$newListAbility = "CWL removed invalid listAbility ID: {0} from {1}"
$ws.Range("C121").Value = $newListAbility
$ws.Range("D121").Value = $newListAbility

# Restore the view state roughly to what was captured in the saved file
# (window scrolled back to the top, selection on D11).
$ws.Range("D11").Select()
